$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "Succinate"
$ws.Range("B16").Value = "Y = 2484862*X + 142079"
$ws.Range("A17").Value = "Fumarate"
$ws.Range("B17").Value = "Y = 5931770*X + 4598311"
$ws.Range("A18").Value = "Citrate"
$ws.Range("B18").Value = "Y = 78356408*X + 10610249"

$newFontRange = $ws.Range("A16:A18")
$newFontRange.Font.Name = "Arial"
$newFontRange.Font.Size = 11
$newFontRange.Font.Color = 1907741
